# Auto-generated edit script: applies numeric cell updates to the
# Kraken_Profits-style leve-profit workbook (per-sheet tables: ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR) to match the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 405.11765
$ws.Range("I15").Value = 405.11765
$ws.Range("K15").Value = 1215.35295
$ws.Range("M15").Value = -1046.35295
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = ""
$ws.Range("H63").Value = 86000
$ws.Range("H66").Value = 86000
$ws.Range("H100").Value = 2834.6667
$ws.Range("I100").Value = 2834.6667
$ws.Range("K100").Value = 2834.6667
$ws.Range("M100").Value = -2293.6667
$ws.Range("H103").Value = 750.4286
$ws.Range("I103").Value = 730.6
$ws.Range("K103").Value = 2191.8
$ws.Range("M103").Value = -1605.8
$ws.Range("H141").Value = 4386.75
$ws.Range("I141").Value = 4386.75
$ws.Range("K141").Value = 13160.25
$ws.Range("M141").Value = -7980.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 99995
$ws.Range("J117").Value = 99995
$ws.Range("L117").Value = 99995
$ws.Range("N117").Value = -109173
$ws.Range("H121").Value = 99995
$ws.Range("J121").Value = 99995
$ws.Range("L121").Value = 99995
$ws.Range("N121").Value = -103489
$ws.Range("H133").Value = 99995
$ws.Range("J133").Value = 99995
$ws.Range("L133").Value = 99995
$ws.Range("N133").Value = -105055
$ws.Range("H140").Value = 47214.5
$ws.Range("J140").Value = 47214.5
$ws.Range("L140").Value = 47214.5
$ws.Range("N140").Value = -57574.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 13333.333
$ws.Range("I96").Value = 13333.333
$ws.Range("K96").Value = 13333.333
$ws.Range("M96").Value = -10587.333
$ws.Range("H107").Value = 3600
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 400
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 1520
$ws.Range("N107").Value = -13840
$ws.Range("H133").Value = 99995
$ws.Range("J133").Value = 99995
$ws.Range("L133").Value = 99995
$ws.Range("N133").Value = -110115

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4191.154
$ws.Range("I31").Value = 4875
$ws.Range("K31").Value = 4875
$ws.Range("M31").Value = -4580
$ws.Range("H34").Value = 4191.154
$ws.Range("I34").Value = 4875
$ws.Range("K34").Value = 4875
$ws.Range("M34").Value = -4673
$ws.Range("H58").Value = 477
$ws.Range("I58").Value = 477
$ws.Range("K58").Value = 477
$ws.Range("M58").Value = -274
$ws.Range("H59").Value = 72000
$ws.Range("J59").Value = 72000
$ws.Range("L59").Value = 72000
$ws.Range("N59").Value = -74290
$ws.Range("H62").Value = 10334
$ws.Range("I62").Value = 13001
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 13001
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -12377
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 10334
$ws.Range("I65").Value = 13001
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 65005
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -61885
$ws.Range("N65").Value = -31240
$ws.Range("H94").Value = 1385.4286
$ws.Range("I94").Value = 1416.3334
$ws.Range("K94").Value = 1416.3334
$ws.Range("M94").Value = -965.3334
$ws.Range("H107").Value = 100000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 100000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 100000
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = -103840
$ws.Range("H132").Value = 4529.6665
$ws.Range("I132").Value = 3946
$ws.Range("K132").Value = 11838
$ws.Range("M132").Value = -9308
$ws.Range("H134").Value = 3346.75
$ws.Range("I134").Value = 3462.6667
$ws.Range("J134").Value = 2999
$ws.Range("K134").Value = 10388.0001
$ws.Range("L134").Value = 8997
$ws.Range("M134").Value = -7853.000100000001
$ws.Range("N134").Value = -14067
$ws.Range("H136").Value = 477
$ws.Range("I136").Value = 477
$ws.Range("K136").Value = 1431
$ws.Range("M136").Value = 1119

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 5000000
$ws.Range("I87").Value = 5000000
$ws.Range("K87").Value = 15000000
$ws.Range("M87").Value = -14998752
$ws.Range("H90").Value = 5000000
$ws.Range("I90").Value = 5000000
$ws.Range("K90").Value = 45000000
$ws.Range("M90").Value = -44993760
$ws.Range("H131").Value = 2315.75
$ws.Range("J131").Value = 2616.5
$ws.Range("L131").Value = 7849.5
$ws.Range("N131").Value = -17929.5
$ws.Range("H138").Value = 2529.875
$ws.Range("I138").Value = 2462.7144
$ws.Range("K138").Value = 7388.1432
$ws.Range("M138").Value = -2248.1432

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
$ws.Range("H132").Value = 2107.6155
$ws.Range("I132").Value = 1639.9
$ws.Range("J132").Value = 3666.6667
$ws.Range("K132").Value = 4919.700000000001
$ws.Range("L132").Value = 11000.0001
$ws.Range("M132").Value = -2389.700000000001
$ws.Range("N132").Value = -16060.0001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5138
$ws.Range("I7").Value = 5093.2
$ws.Range("K7").Value = 5093.2
$ws.Range("M7").Value = -4981.2
$ws.Range("H22").Value = 2000.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 2000.5
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""
$ws.Range("H40").Value = 8100.6665
$ws.Range("I40").Value = 7719.8
$ws.Range("K40").Value = 7719.8
$ws.Range("M40").Value = -7583.8
$ws.Range("H46").Value = 3679.0715
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 3958.9167
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 3958.9167
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -4334.9167
$ws.Range("H76").Value = 20684.625
$ws.Range("J76").Value = 20884.428
$ws.Range("L76").Value = 20884.428
$ws.Range("N76").Value = -21560.428
$ws.Range("H79").Value = 20684.625
$ws.Range("J79").Value = 20884.428
$ws.Range("L79").Value = 20884.428
$ws.Range("N79").Value = -23224.428
$ws.Range("H122").Value = 9500
$ws.Range("I122").Value = 9500
$ws.Range("K122").Value = 28500
$ws.Range("M122").Value = -26050
$ws.Range("H126").Value = 5138
$ws.Range("I126").Value = 5093.2
$ws.Range("K126").Value = 15279.6
$ws.Range("M126").Value = -12809.6
$ws.Range("H132").Value = 8558.799999999999
$ws.Range("I132").Value = 8698.5
$ws.Range("K132").Value = 26095.5
$ws.Range("M132").Value = -23565.5
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 67429
$ws.Range("J46").Value = 67429
$ws.Range("L46").Value = 67429
$ws.Range("N46").Value = -67891
$ws.Range("H131").Value = 99995
$ws.Range("J131").Value = 99995
$ws.Range("L131").Value = 99995
$ws.Range("N131").Value = -110075
$ws.Range("H134").Value = 67429
$ws.Range("J134").Value = 67429
$ws.Range("L134").Value = 202287
$ws.Range("N134").Value = -207357
$ws.Range("H140").Value = 11499.5
$ws.Range("J140").Value = 11499.5
$ws.Range("L140").Value = 11499.5
$ws.Range("N140").Value = -21859.5

